$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.863.17'
$ws.Range('E2').Value = '  +2.66%  '
$ws.Range('D3').Value = '1.666.52'
$ws.Range('D5').Value = '214.64'
$ws.Range('E5').Value = '  +0.12%  '
$ws.Range('E6').Value = '  -0.68%  '
$ws.Range('E8').Value = '  +3.10%  '
$ws.Range('E9').Value = '  +0.32%  '
$ws.Range('D10').Value = '0.0621'
$ws.Range('E10').Value = '  -0.16%  '
$ws.Range('E11').Value = '  -1.43%  '
$ws.Range('D12').Value = '1.901.27'
$ws.Range('E12').Value = '  -0.76%  '
$ws.Range('D13').Value = '1.663.15'
$ws.Range('E13').Value = '  -0.74%  '
$ws.Range('E14').Value = '  -1.50%  '
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').Value = '65.85'
$ws.Range('E16').Value = '  -1.11%  '
$ws.Range('D17').Value = '249.83'
$ws.Range('E17').Value = '  +6.22%  '
$ws.Range('D18').Value = '27.845.93'
$ws.Range('E18').Value = '  +2.71%  '
$ws.Range('D19').Value = '0.0₃0731'
$ws.Range('E19').Value = '  -1.20%  '
$ws.Range('D20').Value = '7.55'
$ws.Range('E20').Value = '  -4.34%  '
$ws.Range('E21').Value = '  -0.02%  '
$ws.Range('D22').Value = '4.47'
$ws.Range('E22').Value = '  -1.52%  '
$ws.Range('E23').Value = '  -1.78%  '
$ws.Range('E24').Value = '  -1.66%  '
$ws.Range('D25').Value = '146.79'
$ws.Range('E25').Value = '  -1.23%  '
$ws.Range('D26').Value = '7.22'
$ws.Range('E26').Value = '  -2.84%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('D29').Value = '0.112'
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('D30').Value = '1.24'
$ws.Range('E30').Value = '  +5.88%  '
$ws.Range('E31').Value = '  -0.46%  '
$ws.Range('D32').Value = '3.34'
$ws.Range('E32').Value = '  -0.55%  '
$ws.Range('E33').Value = '  -3.11%  '
$ws.Range('D34').Value = '1.412.80'
$ws.Range('E34').Value = '  -8.29%  '
$ws.Range('E35').Value = '  -5.81%  '
$ws.Range('E36').Value = '  -0.06%  '
$ws.Range('D37').Value = '0.930'
$ws.Range('E37').Value = '  -0.97%  '
$ws.Range('D38').Value = '0.579'
$ws.Range('E38').Value = '  -4.66%  '
$ws.Range('E39').Value = '  -1.86%  '
$ws.Range('E40').Value = '  -2.54%  '
$ws.Range('D41').Value = '69.37'
$ws.Range('E41').Value = '  -0.10%  '
$ws.Range('E42').Value = '  -0.01%  '
$ws.Range('D43').Value = '5.44'
$ws.Range('E43').Value = '  -5.99%  '
$ws.Range('E44').Value = '  -1.12%  '
$ws.Range('B45').Value = 'RocketPoolETH'
$ws.Range('C45').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D45').Value = '1.809.74'
$ws.Range('E45').Value = '  -0.85%  '
$ws.Range('B46').Value = 'TrustWalletToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D46').Value = '0.790'
$ws.Range('E46').Value = '  +1.22%  '
$ws.Range('E47').Value = '  +4.95%  '
$ws.Range('D48').Value = '88.50'
$ws.Range('E48').Value = '  -1.64%  '
$ws.Range('D49').Value = '0.0₆0110'
$ws.Range('E49').Value = '  -0.81%  '
$ws.Range('E50').Value = '  -2.68%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').Value = '0.0510'
$ws.Range('E51').Value = '  -0.16%  '
